$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4 share the same duty date (column A) and shift (column B). Clear
# the now-redundant values in rows 3-4 before merging those cells with row 2.
$ws.Range("A3:A4").ClearContents()
$ws.Range("B3:B4").ClearContents()

$ws.Range("B2:B4").Merge()
$ws.Range("A2:A4").Merge()

# The merged date cell also picks up horizontal centering (it already had
# vertical centering).
$ws.Range("A2:A4").HorizontalAlignment = -4108

# Move the active selection from D17 to D12.
$ws.Range("D12").Select()
